$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): swap C and pi columns ---
$ws.Range("C1").Value = "$\pi$"
$ws.Range("D1").Value = "C"

# --- Row 2 stays "A Lag", update coefficients ---
$ws.Range("B2").Value = "-0.412***"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "0.202"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1.004"
$ws.Range("D2").Style = "Normal"

# --- Row 3 becomes "$\pi$ Lag" (was "C Lag"), update coefficients ---
$ws.Range("A3").Value = "$\pi$ Lag"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "-0.071"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = "-0.636***"
$ws.Range("D3").Value = "1.118**"

# --- Row 4 becomes "C Lag" (was "$\pi$ Lag"), update coefficients ---
$ws.Range("A4").Value = "C Lag"
$ws.Range("B4").Value = "-0.052***"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "-0.013"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "-0.698***"

# --- Remove the old "Constant" and "r2_adj" rows entirely ---
$ws.Rows("5:6").Delete()
